$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New labels for column A, rows 63-70 (school2_4yr/school2_ad/school2_cc removed,
# "school2" renamed to "school2_type"; B/C values for these rows are left untouched).
$newLabels = @(
    "school2_type",
    "sexual",
    "sib_freq",
    "stig_pcv_2",
    "stig_pcv_3",
    "talk",
    "ther_vis",
    "wcs_tot"
)

for ($i = 0; $i -lt $newLabels.Length; $i++) {
    $row = 63 + $i
    $ws.Range("A$row").Value = $newLabels[$i]
}

# Rows 71-73 no longer exist; delete them entirely (shifting cells up, which here
# just removes their now-stale content since nothing is below them).
$ws.Range("A71:C73").Delete() | Out-Null
